# Update drs_data worksheet: append 4 new DRS review rows (Match 26, LSG vs DC)
# and move the active selection/top-left cell as in the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append starting at row 87.
$newRows = @(
    @("26","LSG","DC","1","LSG","DC","3","LSG","R Pandit","RP","Wicket","Out","Out","Q de Kock","KK Ahmed","Unsuccessful","Yes"),
    @("26","LSG","DC","1","LSG","DC","4","DC","YC Barde","YCB","Wide","Called","Called","D Padikkal","I Sharma","Unsuccessful","No"),
    @("26","LSG","DC","1","LSG","DC","10","DC","YC Barde","YCB","Wicket","Not Out","Out","KL Rahul","Kuldeep Yadav","Successful",$null),
    @("26","LSG","DC","1","LSG","DC","13","DC","R Pandit","RP","Wicket","Not Out","Not Out","A Badoni","Mukesh Kumar","Unsuccessful","No")
)

$startRow = 87
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $val = $rowValues[$c - 1]
        if ($null -eq $val) {
            # Leave the cell completely untouched/blank (no cell element emitted).
            continue
        }
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $val
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4108
    }
}

# Columns A, D, G hold numeric values; reset them to numbers (not text) to match source data.
$numericCols = @(1, 4, 7)
foreach ($i in 0..($newRows.Count - 1)) {
    $r = $startRow + $i
    foreach ($c in $numericCols) {
        $v = $newRows[$i][$c - 1]
        $ws.Cells.Item($r, $c).Value = [double]$v
    }
}

# Update the view: scroll the window so column A is the left-most visible
# column again (top-left cell moves from C62 to A62), and move the active
# selection to I89 (previously N85).
$ws.Application.ActiveWindow.ScrollRow = 62
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I89").Select()
